$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, shifting existing rows 81-83 down to 82-84
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with the new weekly data point
$ws.Range("A81").Value = 11
$ws.Range("B81").Value = "Vega Monumental Concepción"
$ws.Range("C81").Value = "Bíobío"
$ws.Range("D81").Value = 44509
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = 100112043
$ws.Range("G81").Value = "Pepino ensalada"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 100
$ws.Range("K81").Value = 6000
$ws.Range("L81").Value = 6500
$ws.Range("M81").Value = 6250
$ws.Range("N81").Value = "$/caja 60 unidades"
$ws.Range("O81").Value = "Región de Arica y Parinacota"
$ws.Range("P81").Value = 104
$ws.Range("Q81").Value = 60
$ws.Range("R81").Value = "Hortaliza"
